$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last status check" timestamp in F1
$ws.Range("F1").Value = "Last status check on: 31.01.2022 01:15"

# Update row 5 (Makro) prices - values swapped
$ws.Range("B5").Value = 34.9
$ws.Range("C5").Value = 34.5

# Delta now stored as text with explicit "+" sign
$ws.Range("D5").Value = "'+0.4"
$ws.Range("D5").ClearFormats()

# Old-date column now stored as plain text timestamp (no date style)
$ws.Range("E5").Value = "2022-01-31 01:15:09"
$ws.Range("E5").ClearFormats()
